$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: waitfortext step now waits for "Mattresses" text (was the ecommerce page title)
$ws.Range("D3").Value = "Mattresses"

# Row 4: turn the old "assert" step into a "click" step on the "Mattresses" link
$ws.Range("D4").Value = "Mattresses"
$ws.Range("C4").Value = "click"
$ws.Range("E4").Value = "a"
$ws.Range("F4").ClearContents()
$ws.Range("H4").Value = 2000

# Row 2: goto target changes from the scraping-course demo site to nectarsleep.com
$ws.Range("D2").Value = "https://www.nectarsleep.com"

# Make the goto URL a real hyperlink
$ws.Hyperlinks.Add($ws.Range("D2"), "https://www.nectarsleep.com")

# Hyperlinks.Add() re-applies a fresh "Hyperlink" style variant to the cell;
# restore the original Hyperlink-style formatting (matching E2) that D2 already had
$ws.Range("E2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
